# Generate Report for Handback
#
# Refresh the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the f15435dc-... handback
# record across the Overview, zh-cn and de-de sheets, reflecting a new
# handback report run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# f15435dc-7905-498a-9033-17377cc1815a.md row (row 3).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-31 21:02:09"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the f15435dc-7905-498a-9033-17377cc1815a row (row 3).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-31 21:01:57"
$wsZhCn.Range("K3").Value = "2016-08-31 21:02:35"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the f15435dc-7905-498a-9033-17377cc1815a row (row 3).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-31 21:02:09"
$wsDeDe.Range("K3").Value = "2016-08-31 21:02:43"
